$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D3").Value = 'Background/objectives
id="Par1">To evaluate the prevalence of malnutrition and its related factors in elderly patients with COVID-19 in Wuhan, China.


Subjects/methods
id="Par2">In a cross-sectional study, we evaluated the nutritional status of elderly inpatients with COVID-19 using the Mini Nutritional Assessment (MNA).

 Based on MNA scores, patients were divided into non-malnutrition group (MNA ≥ 24), the group with risk of malnutrition (MNA 17–23.5) and malnutrition group (MNA score &lt; 17).

 Regression analysis was conducted to screen for risk factors for malnutrition.


Results
id="Par3">A total of 182 patients were included in the study, of which 27.5% were in the group with malnutrition risk and 52.7% were in the malnutrition group.

 There were statistical differences in the incidence of comorbid diabetes mellitus, body mass index (BMI), calf circumference,albumin, hemoglobin, and lymphocyte counts among the three groups.

 Further regression analysis suggested that combined diabetes, low calf circumference, and low albumin were independent risk factors for malnutrition.


Conclusions
id="Par4">The prevalence of malnutrition in elderly patients with COVID-19 was high, and nutritional support should be strengthened during treatment, especially for those with diabetes mellitus, low calf circumference, or low albumin.


'
$ws.Range("E3").Value = '[Tao%Li%NULL%1,    Yalan%Zhang%NULL%2,    Yalan%Zhang%NULL%0,    Cheng%Gong%NULL%1,    Jing%Wang%NULL%0,    Bao%Liu%NULL%1,    Li%Shi%NULL%1,    Jun%Duan%junjununun@163.com%1]'
$ws.Range("D5").Value = 'Objectives
id="Par1">The aim of this study was to assess the nutritional risks among older patients with COVID-19 and their associated clinical outcomes using four nutritional risk screening (NRS) tools: Nutrition Risk Screening 2002 (NRS 2002), Malnutrition Universal Screening Tool (MUST), Mini Nutrition Assessment Shortcut (MNA-sf), and Nutrition Risk Index (NRI).


Methods
id="Par2">We retrospectively analyzed the data of patients with COVID-19 older than 65 years who were treated in our hospital from January 28, 2020 to March 5, 2020, and explored the relationship between nutritional risk and clinical outcomes.


Results
id="Par3">A total of 141 patients with COVID-19 (46 common COVID-19, 73 severe COVID-19, and 22 extremely severe COVID-19) were enrolled in the study.

 NRS 2002 identified 85.8% of patients as having risk, with being identified 41.1% by MUST, 77.3% by MNA-sf, and 71.6% by NRI.

 The agreement strength was moderate between NRS 2002 and MNA-sf, NRI, fair between MUST and MNA-sf, NRI, fair between MNA-sf and NRI, poor between NRS 2002 and MUST (P &lt; 0.01).

 After adjustment for confounding factors in multivariate regression analysis, patients in the risk group had significantly longer LOS, higher hospital expenses (except MNA-sf), poor appetite, heavier disease severity, and more weight change(kg) than normal patients by using NRS 2002, MNA-sf, and NRI(P &lt; 0.05).


Conclusions
id="Par4">The NRS 2002, MNA-sf, and NRI are useful and practical tools with respect to screening for patients with COVID-19 who are at nutritional risk, as well as in need of additional nutritional intervention.


'
$ws.Range("E5").Value = '[Gaoli%Liu%NULL%1,    Shaowen%Zhang%NULL%2,    Shaowen%Zhang%NULL%0,    Zhangfan%Mao%NULL%1,    Weixing%Wang%13392186@qq.com%1,    Haifeng%Hu%NULL%1]'
